$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column H (copy formatting from G1 so it matches the other headers)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Save values for rows 2-28 (matches "Save" flag per row)
$saveValues = @(1,0,0,0,1,0,0,1,1,0,0,0,0,0,1,0,1,0,1,0,0,1,0,0,0,1,0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
